# Insert a new weekly price-report row for "Rabanito" (Vega Modelo de Temuco)
# at row 36, pushing the existing rows 36-95 down to 37-96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 36..95 down to 37..96 and create a fresh blank row 36
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new record
$ws.Cells.Item(36, 1).Value  = 10
$ws.Cells.Item(36, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(36, 3).Value  = "La Araucanía"
$ws.Cells.Item(36, 4).Value  = 44952
$ws.Cells.Item(36, 5).Value  = 9
$ws.Cells.Item(36, 6).Value  = 300000001
$ws.Cells.Item(36, 7).Value  = "Rabanito"
$ws.Cells.Item(36, 8).Value  = "Sin especificar"
$ws.Cells.Item(36, 9).Value  = "Primera"
$ws.Cells.Item(36, 10).Value = 65
$ws.Cells.Item(36, 11).Value = 8000
$ws.Cells.Item(36, 12).Value = 8000
$ws.Cells.Item(36, 13).Value = 8000
$ws.Cells.Item(36, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(36, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(36, 16).Value = 667
$ws.Cells.Item(36, 17).Value = 12
$ws.Cells.Item(36, 18).Value = "Hortaliza"
